$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Fill in the self-evaluation level values in column B
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 5
$ws.Range("B9").Value = 5
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 3
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 5
$ws.Range("B15").Value = 5
$ws.Range("B16").Value = 5
$ws.Range("B17").Value = 3
$ws.Range("B18").Value = 2
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 1
$ws.Range("B22").Value = 3
$ws.Range("B23").Value = 3
$ws.Range("B24").Value = 4
$ws.Range("B25").Value = 3
$ws.Range("B26").Value = 4
$ws.Range("B27").Value = 4
$ws.Range("B28").Value = 5
$ws.Range("B29").Value = 3
$ws.Range("B30").Value = 1
$ws.Range("B32").Value = 5
$ws.Range("B33").Value = 5
$ws.Range("B34").Value = 5
$ws.Range("B35").Value = 4
$ws.Range("B36").Value = 4
$ws.Range("B39").Value = 3
$ws.Range("B40").Value = 4
$ws.Range("B41").Value = 4
$ws.Range("B42").Value = 3
$ws.Range("B43").Value = 4
$ws.Range("B45").Value = 5
$ws.Range("B46").Value = 5
$ws.Range("B47").Value = 4

# Update the sheet view: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B48").Select()
